$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 - Mendonça et al. 117th Congress BERTopic article
$ws.Range("C6").Value = "Topic Mining, BERTopic, 117th Congress, Twitter, short-text data"
$ws.Range("E6").Value = 5

# Row 7 - Blei & Lafferty Dynamic Topic Models article
$ws.Range("C7").Value = "Latent Dirichlet Allocation (LDA), Dynamic Topic Models (DTM)"
$ws.Range("E7").Value = 6

# Row 8 - Wang & McCallum Topics over Time article
$ws.Range("C8").Value = "Graphical Models, Temporal Analysis, Topic Modeling"
$ws.Range("E8").Value = 7

# Row 9 - Teh et al. Hierarchical Dirichlet Processes article
$ws.Range("C9").Value = "Clustering, Hierarchical model, Markov chain Monte Carlo, Mixture model, Nonparametric Bayesian statistics"
$ws.Range("E9").Value = 8

# Add hyperlinks to D6:D9 matching their displayed URL text
$ws.Hyperlinks.Add($ws.Range("D6"), $ws.Range("D6").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), $ws.Range("D7").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), $ws.Range("D8").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), $ws.Range("D9").Text) | Out-Null

# Update selection to reflect the final active cell used while editing
$ws.Range("C10").Select() | Out-Null
